# Update countries & provincias Spain
# - Refresh "last updated" timestamp
# - Update case counts for several countries (Estados Unidos, Reino Unido, Rumania, Sri Lanka)
# - Re-sort "Santa Lucia" into its alphabetically correct place (before Namibia) with
#   refreshed data, shifting Namibia / Curazao / Dominica / San Cristobal y Nieves /
#   San Vicente y las Granadinas down by one row each

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" banner text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 19:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1046426
$ws.Range("C4").Value = 10661
$ws.Range("D4").Value = 144352
$ws.Range("E4").Value = 841962
$ws.Range("G4").Value = 846
$ws.Range("H4").Value = 60112

# --- Row 8: Reino Unido ---
$ws.Range("B8").Value = 165221
$ws.Range("C8").Value = 4076
$ws.Range("E8").Value = 138780
$ws.Range("G8").Value = 4419
$ws.Range("H8").Value = 26097

# --- Row 36: Rumania ---
$ws.Range("E36").Value = 7721
$ws.Range("G36").Value = 25
$ws.Range("H36").Value = 688

# --- Row 104: Sri Lanka ---
$ws.Range("B104").Value = 630
$ws.Range("C104").Value = 11
$ws.Range("E104").Value = 487

# --- Move "Santa Lucia" up to its sorted position (currently row 195) so it sits
#     right before "Namibia" (currently row 190), refreshing its data in the process.
#     Deleting the old row first (and inserting the new one afterwards) makes the
#     shift land on the correct final row numbers. ---
$ws.Rows(195).Delete()
$ws.Rows(190).Insert()

$ws.Range("A190").Value = "Santa Lucia"
$ws.Range("B190").Value = 17
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 15
$ws.Range("E190").Value = 2
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0
